$d = $word.ActiveDocument

# Fix typo in placeholder: {week_sarcharit_1} -> {week_shacharit_1}
$d.Content.Find.Execute("{week_sarcharit_1}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{week_shacharit_1}", 2)
